$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.619.13"
$ws.Range("E2").Value = "  +0.41%  "

# Row 3
$ws.Range("D3").Value = "3.512.04"
$ws.Range("E3").Value = "  -0.04%  "

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.10%  "

# Row 5
$ws.Range("D5").Value = "'609.56"
$ws.Range("E5").Value = "  -0.13%  "

# Row 6
$ws.Range("D6").Value = "'152.50"
$ws.Range("E6").Value = "  +1.43%  "

# Row 7
$ws.Range("D7").Value = "3.510.88"
$ws.Range("E7").Value = "  -0.03%  "

# Row 8
$ws.Range("E8").Value = "  -0.04%  "

# Row 9
$ws.Range("E9").Value = "  +1.35%  "

# Row 10
$ws.Range("E10").Value = "  +3.27%  "

# Row 11
$ws.Range("D11").Value = "'7.59"
$ws.Range("E11").Value = "  +7.80%  "

# Row 12
$ws.Range("D12").Value = "'0.433"
$ws.Range("E12").Value = "  +1.76%  "

# Row 13
$ws.Range("D13").Value = "'32.61"
$ws.Range("E13").Value = "  +2.58%  "

# Row 14
$ws.Range("E14").Value = "  -1.25%  "

# Row 15
$ws.Range("D15").Value = "4.103.17"
$ws.Range("E15").Value = "  -0.12%  "

# Row 16
$ws.Range("D16").Value = "3.504.85"
$ws.Range("E16").Value = "  -0.34%  "

# Row 17
$ws.Range("D17").Value = "67.469.07"
$ws.Range("E17").Value = "  +0.17%  "

# Row 18
$ws.Range("E18").Value = "  +0.09%  "

# Row 19
$ws.Range("D19").Value = "'6.57"
$ws.Range("E19").Value = "  +2.63%  "

# Row 20
$ws.Range("D20").Value = "'15.61"
$ws.Range("E20").Value = "  +2.45%  "

# Row 21
$ws.Range("D21").Value = "'9.87"
$ws.Range("E21").Value = "  +6.72%  "

# Row 22
$ws.Range("D22").Value = "'447.20"
$ws.Range("E22").Value = "  +0.98%  "

# Row 23
$ws.Range("E23").Value = "  +1.30%  "

# Row 24
$ws.Range("D24").Value = "'78.19"
$ws.Range("E24").Value = "  +1.26%  "

# Row 25
$ws.Range("D25").Value = "3.650.19"
$ws.Range("E25").Value = "  -0.13%  "

# Row 26
$ws.Range("E26").Value = "  +0.07%  "

# Row 27
$ws.Range("E27").Value = "  -1.35%  "

# Row 28
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "'8.81"
$ws.Range("E28").Value = "  +5.32%  "

# Row 29
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value = "'10.12"
$ws.Range("E29").Value = "  -1.31%  "

# Row 30
$ws.Range("D30").Value = "'2.51"
$ws.Range("E30").Value = "  +0.61%  "

# Row 31
$ws.Range("E31").Value = "  +7.45%  "

# Row 32
$ws.Range("E32").Value = "  +4.59%  "

# Row 33
$ws.Range("E33").Value = "  +0.03%  "

# Row 34
$ws.Range("D34").Value = "'25.73"
$ws.Range("E34").Value = "  -0.09%  "

# Row 35
$ws.Range("D35").Value = "'6.19"
$ws.Range("E35").Value = "  +0.84%  "

# Row 36
$ws.Range("E36").Value = "  +1.89%  "

# Row 37
$ws.Range("D37").Value = "3.505.51"
$ws.Range("E37").Value = "  +0.01%  "

# Row 38
$ws.Range("E38").Value = "  +0.60%  "

# Row 39
$ws.Range("E39").Value = "  +0.07%  "

# Row 40
$ws.Range("D40").Value = "'2.32"
$ws.Range("E40").Value = "  +7.40%  "

# Row 41
$ws.Range("E41").Value = "  -0.11%  "

# Row 42
$ws.Range("E42").Value = "  +2.98%  "

# Row 43
$ws.Range("D43").Value = "'173.13"
$ws.Range("E43").Value = "  -2.77%  "

# Row 44
$ws.Range("D44").Value = "'5.48"
$ws.Range("E44").Value = "  +0.77%  "

# Row 45
$ws.Range("D45").Value = "'30.27"
$ws.Range("E45").Value = "  +10.03%  "

# Row 46
$ws.Range("D46").Value = "'0.885"
$ws.Range("E46").Value = "  +0.56%  "

# Row 47
$ws.Range("D47").Value = "'46.66"
$ws.Range("E47").Value = "  +2.36%  "

# Row 48
$ws.Range("E48").Value = "  +3.77%  "

# Row 49
$ws.Range("E49").Value = "  -2.72%  "

# Row 50
$ws.Range("D50").Value = "'7.66"
$ws.Range("E50").Value = "  +1.14%  "

# Row 51
$ws.Range("D51").Value = "'0.996"
$ws.Range("E51").Value = "  -0.05%  "
